# TIRANDO AS MENSAGENS ANTERIORES ANTES DE GERAR O QR CODE
#
# Fill in the payment id (idPagamento) for three previously-pending rows,
# and append three fresh "Vitor Ito" rows (new numbers-selected entries,
# still awaiting payment, so their idPagamento stays blank) below the
# existing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    # Force the cell to be stored as text even when the value looks numeric
    # (phone numbers / payment ids), mirroring how Excel keeps a cell's
    # content as text when it was typed into a Text-formatted cell - then
    # drop the number format again so the cell keeps the workbook's normal
    # (General) styling.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- Fill in idPagamento for the three rows that were waiting on it ---
Set-TextValue $ws.Range("D24") "76956334770"
Set-TextValue $ws.Range("D25") "76956534748"
Set-TextValue $ws.Range("D26") "76956509396"

# --- Append the three new rows (30-32) ---
$newRows = @(
    @{ Row = 30; Phone = "11988776655" },
    @{ Row = 31; Phone = "11977665544" },
    @{ Row = 32; Phone = "11887766554" }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    $ws.Range("A$r").Value = "Vitor Ito"
    $ws.Range("B$r").Value = 1578424633
    Set-TextValue $ws.Range("C$r") $entry.Phone
    # D$r (idPagamento) intentionally left blank - no payment id yet.
    $ws.Range("E$r").Value = 1
    $ws.Range("F$r").Value = 2
    $ws.Range("G$r").Value = 3
    $ws.Range("H$r").Value = 4
    $ws.Range("I$r").Value = 5
    $ws.Range("J$r").Value = 6
    $ws.Range("K$r").Value = 7
    $ws.Range("L$r").Value = 8
    $ws.Range("M$r").Value = 9
    $ws.Range("N$r").Value = 10
    $ws.Range("O$r").Value = "Não"
}
